$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string "Baza podataka" -> "Analiza sadržaja"
# (this cell currently holds " Baza podataka")
$ws.Range("D3").Value = "Analiza sadržaja"

# Update current selection to D3 only (previously A3:D9 with active cell A3)
$ws.Range("D3").Select()
